# Add a new worksheet "studyDesignOE" after "studyDesignPopulations" and
# populate it with the study objectives / endpoints table.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new sheet at the end of the workbook ----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "studyDesignOE"

# --- 2. Column widths ----------------------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 43.166666666666664   # -> 44
$newSheet.Columns.Item(2).ColumnWidth = 17.998697916666668   # -> ~18.83203125
$newSheet.Columns.Item(3).ColumnWidth = 58.666666666666664   # -> 59.5
$newSheet.Columns.Item(4).ColumnWidth = 30.666666666666668   # -> 31.5
$newSheet.Columns.Item(5).ColumnWidth = 30.666666666666668   # -> 31.5

# --- 3. Build the two cell styles used on this sheet ---------------------
# Style A: bold header font, wrap text, left/top aligned (header row).
$seedHeader = $newSheet.Range("Z101")
$seedHeader.Font.Bold = $true
$seedHeader.WrapText = $true
$seedHeader.HorizontalAlignment = -4131
$seedHeader.VerticalAlignment = -4160
$seedHeader.Copy()
$newSheet.Range("A1:E1").PasteSpecial(-4122)
$seedHeader.Clear()

# Style B: regular font, wrap text, left/top aligned (body rows).
$seedBody = $newSheet.Range("Z100")
$seedBody.WrapText = $true
$seedBody.HorizontalAlignment = -4131
$seedBody.VerticalAlignment = -4160
$seedBody.Copy()
$newSheet.Range("A2:E36").PasteSpecial(-4122)
$seedBody.Clear()

# --- 4. Row heights --------------------------------------------------------
$newSheet.Rows.Item(1).RowHeight = 17
$newSheet.Rows.Item(2).RowHeight = 68
$newSheet.Rows.Item(3).RowHeight = 68
$newSheet.Rows.Item(4).RowHeight = 34
$newSheet.Rows.Item(5).RowHeight = 17
$newSheet.Rows.Item(6).RowHeight = 17
$newSheet.Rows.Item(7).RowHeight = 17
$newSheet.Rows.Item(8).RowHeight = 17
$newSheet.Rows.Item(9).RowHeight = 85
$newSheet.Rows.Item(10).RowHeight = 17
$newSheet.Rows.Item(11).RowHeight = 51
$newSheet.Rows.Item(12).RowHeight = 85
$newSheet.Rows.Item(13).RowHeight = 17

# --- 5. Cell values (row-major order so shared strings line up) ----------
$newSheet.Range("A1").Value = 'objectiveDescription'
$newSheet.Range("B1").Value = 'objectiveLevel'
$newSheet.Range("C1").Value = 'endpointDescription'
$newSheet.Range("D1").Value = 'endpointPurposeDescription'
$newSheet.Range("E1").Value = 'endpointLevel'

$newSheet.Range("A2").Value = 'The primary efficacy objective for this study is to evaluate the efficacy of TCZ compared with placebo in combination with SOC for the treatment of severe COVID-19 pneumonia'
$newSheet.Range("B2").Value = 'Study Primary Objective'
$newSheet.Range("C2").Value = 'Clinical status assessed using a 7-category ordinal scale at Day 28'
$newSheet.Range("E2").Value = 'Primary Endpoint'

$newSheet.Range("A3").Value = 'The secondary efficacy objective for this study is to evaluate the efficacy of TCZ compared with placebo in combination with SOC for the treatment of severe COVID-19 pneumonia'
$newSheet.Range("B3").Value = 'Study Secondary Objective'
$newSheet.Range("C3").Value = 'Time to clinical improvement (TTCI) defined as a National Early Warning Score 2 (NEWS2) of <=2 maintained for 24 hours'
$newSheet.Range("E3").Value = 'Secondary Enpoint'

$newSheet.Range("C4").Value = 'Time to improvement of at least 2 categories relative to baseline on a 7-category ordinal scale of clinical status'
$newSheet.Range("E4").Value = 'Secondary Enpoint'

$newSheet.Range("C5").Value = 'Incidence of mechanical ventilation'
$newSheet.Range("E5").Value = 'Secondary Enpoint'

$newSheet.Range("C6").Value = 'Ventilator-free days to Day 28'
$newSheet.Range("E6").Value = 'Secondary Enpoint'

$newSheet.Range("C7").Value = 'Incidence of intensive care unit (ICU) stay'
$newSheet.Range("E7").Value = 'Secondary Enpoint'

$newSheet.Range("C8").Value = 'Duration of ICU stay'
$newSheet.Range("E8").Value = 'Secondary Enpoint'

$newSheet.Range("C9").Value = 'Time to clinical failure, defined as the time to death, mechanical ventilation, ICU admission, or withdrawal (whichever occurs first). For patients entering the study already in ICU or on mechanical ventilation, clinical failure is defined as a one-category worsening on the ordinal scale, withdrawal or death.'
$newSheet.Range("E9").Value = 'Secondary Enpoint'

$newSheet.Range("C10").Value = 'Mortality rate at Days 7, 14, 21, 28, and 60'
$newSheet.Range("E10").Value = 'Secondary Enpoint'

$newSheet.Range("C11").Value = 'Time to hospital discharge or “ready for discharge” (as evidenced by normal body temperature and respiratory rate, and stable oxygen saturation on ambient air or <= 2L supplemental oxygen)'
$newSheet.Range("E11").Value = 'Secondary Enpoint'

$newSheet.Range("C12").Value = 'Time to recovery, defined as discharged or “ready for discharge” (as evidenced by normal body temperature and respiratory rate, and stable oxygen saturation on ambient air or <= 2L supplemental oxygen); OR, in a non-ICU hospital ward (or “ready for hospital ward”) not requiring supplemental oxygen'
$newSheet.Range("E12").Value = 'Secondary Enpoint'

$newSheet.Range("C13").Value = 'Duration of supplemental oxygen'
$newSheet.Range("E13").Value = 'Secondary Enpoint'

# --- 6. Selection / active sheet -----------------------------------------
$newSheet.Range("A8").Select()
